$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.1343
$ws.Range("B10").Value = 8.7592
$ws.Range("B12").Value = 6.355799999999993
$ws.Range("C13").Value = -12.05999999999999
$ws.Range("B18").Value = 5.369000000000001
$ws.Range("E20").Value = 13.17379999999999

$wb.Save()
